$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.056197333333333
$ws.Range("H2").Value = 3.168592
$ws.Range("I2").Value = 0.01247237710445079
$ws.Range("J2").Value = 0.01398563433468744
$ws.Range("M2").Value = 2.035934
$ws.Range("N2").Value = 6.107802
$ws.Range("O2").Value = 0.03126880699186227
$ws.Range("P2").Value = 0.04430738339814538
$ws.Range("Q2").Value = 2.150348061642666
$ws.Range("R2").Value = 19.353132554784
$ws.Range("S2").Value = 0.0003899963524087938
$ws.Range("T2").Value = 0.0006196668625332625

$ws.Range("G3").Value = 1.056197333333333
$ws.Range("H3").Value = 3.168592
$ws.Range("I3").Value = 0.01247237710445079
$ws.Range("J3").Value = 0.01398563433468744
$ws.Range("O3").Value = 0.0582841555718936
$ws.Range("P3").Value = 0.08258768643246023
$ws.Range("Q3").Value = 4.008186848673778
$ws.Range("R3").Value = 36.073681638064
$ws.Range("S3").Value = 0.0007269419675071338
$ws.Range("T3").Value = 0.001155041182992216

$ws.Range("G4").Value = 1.056197333333333
$ws.Range("H4").Value = 3.168592
$ws.Range("I4").Value = 0.01247237710445079
$ws.Range("J4").Value = 0.01398563433468744
$ws.Range("M4").Value = 0.451002
$ws.Range("N4").Value = 1.353006
$ws.Range("O4").Value = 0.006926695310822388
$ws.Range("P4").Value = 0.009815012926416261
$ws.Range("Q4").Value = 0.4763471097279999
$ws.Range("R4").Value = 4.287123987552
$ws.Range("S4").Value = 0.00008639235600420781
$ws.Range("T4").Value = 0.0001372691817790884

$ws.Range("G5").Value = 1.056197333333333
$ws.Range("H5").Value = 3.168592
$ws.Range("I5").Value = 0.01247237710445079
$ws.Range("J5").Value = 0.01398563433468744
$ws.Range("M5").Value = 57.4814495
$ws.Range("N5").Value = 114.962899
$ws.Range("O5").Value = 0.8828264325012393
$ws.Range("P5").Value = 0.8339669888701803
$ws.Range("Q5").Value = 60.71175367803465
$ws.Range("R5").Value = 364.270522068208
$ws.Range("S5").Value = 0.01101094418393243
$ws.Range("T5").Value = 0.01166355735353869

$ws.Range("G6").Value = 1.056197333333333
$ws.Range("H6").Value = 3.168592
$ws.Range("I6").Value = 0.01247237710445079
$ws.Range("J6").Value = 0.01398563433468744
$ws.Range("M6").Value = 1.347395
$ws.Range("N6").Value = 4.042185
$ws.Range("O6").Value = 0.02069390962418245
$ws.Range("P6").Value = 0.02932292837279799
$ws.Range("Q6").Value = 1.423115005946666
$ws.Range("R6").Value = 12.80803505352
$ws.Range("S6").Value = 0.000258102244598227
$ws.Range("T6").Value = 0.0004100997538441842

$ws.Range("I7").Value = 0.6620593097549599
$ws.Range("J7").Value = 0.7423861014276285
$ws.Range("M7").Value = 2.035934
$ws.Range("N7").Value = 6.107802
$ws.Range("O7").Value = 0.03126880699186227
$ws.Range("P7").Value = 0.04430738339814538
$ws.Range("Q7").Value = 114.1448772356333
$ws.Range("R7").Value = 1027.3038951207
$ws.Range("S7").Value = 0.0207018047738934
$ws.Range("T7").Value = 0.03289318562540838

$ws.Range("I8").Value = 0.6620593097549599
$ws.Range("J8").Value = 0.7423861014276285
$ws.Range("O8").Value = 0.0582841555718936
$ws.Range("P8").Value = 0.08258768643246023
$ws.Range("S8").Value = 0.03858756780757858
$ws.Range("T8").Value = 0.0613119505565216

$ws.Range("I9").Value = 0.6620593097549599
$ws.Range("J9").Value = 0.7423861014276285
$ws.Range("M9").Value = 0.451002
$ws.Range("N9").Value = 1.353006
$ws.Range("O9").Value = 0.006926695310822388
$ws.Range("P9").Value = 0.009815012926416261
$ws.Range("Q9").Value = 25.2854797469
$ws.Range("R9").Value = 227.5693177221
$ws.Range("S9").Value = 0.004585883116365987
$ws.Range("T9").Value = 0.007286529181903948

$ws.Range("I10").Value = 0.6620593097549599
$ws.Range("J10").Value = 0.7423861014276285
$ws.Range("M10").Value = 57.4814495
$ws.Range("N10").Value = 114.962899
$ws.Range("O10").Value = 0.8828264325012393
$ws.Range("P10").Value = 0.8339669888701803
$ws.Range("Q10").Value = 3222.704172386609
$ws.Range("R10").Value = 19336.22503431965
$ws.Range("S10").Value = 0.5844834585352041
$ws.Range("T10").Value = 0.6191255015866716

$ws.Range("I11").Value = 0.6620593097549599
$ws.Range("J11").Value = 0.7423861014276285
$ws.Range("M11").Value = 1.347395
$ws.Range("N11").Value = 4.042185
$ws.Range("O11").Value = 0.02069390962418245
$ws.Range("P11").Value = 0.02932292837279799
$ws.Range("Q11").Value = 75.54185787108334
$ws.Range("R11").Value = 679.8767208397501
$ws.Range("S11").Value = 0.01370059552191775
$ws.Range("T11").Value = 0.02176893447712309

$ws.Range("G12").Value = 0.04559766666666667
$ws.Range("H12").Value = 0.136793
$ws.Range("I12").Value = 0.0005384517417354892
$ws.Range("J12").Value = 0.0006037813885615125
$ws.Range("M12").Value = 2.035934
$ws.Range("N12").Value = 6.107802
$ws.Range("O12").Value = 0.03126880699186227
$ws.Range("P12").Value = 0.04430738339814538
$ws.Range("Q12").Value = 0.09283383988733332
$ws.Range("R12").Value = 0.835504558986
$ws.Range("S12").Value = 0.00001683674358675908
$ws.Range("T12").Value = 0.00002675197347165952

$ws.Range("G13").Value = 0.04559766666666667
$ws.Range("H13").Value = 0.136793
$ws.Range("I13").Value = 0.0005384517417354892
$ws.Range("J13").Value = 0.0006037813885615125
$ws.Range("O13").Value = 0.0582841555718936
$ws.Range("P13").Value = 0.08258768643246023
$ws.Range("Q13").Value = 0.1730396035812222
$ws.Range("R13").Value = 1.557356432231
$ws.Range("S13").Value = 0.00003138320508326833
$ws.Range("T13").Value = 0.00004986490799227362

$ws.Range("G14").Value = 0.04559766666666667
$ws.Range("H14").Value = 0.136793
$ws.Range("I14").Value = 0.0005384517417354892
$ws.Range("J14").Value = 0.0006037813885615125
$ws.Range("M14").Value = 0.451002
$ws.Range("N14").Value = 1.353006
$ws.Range("O14").Value = 0.006926695310822388
$ws.Range("P14").Value = 0.009815012926416261
$ws.Range("Q14").Value = 0.020564638862
$ws.Range("R14").Value = 0.185081749758
$ws.Range("S14").Value = 0.000003729691154583361
$ws.Range("T14").Value = 0.000005926122133460804

$ws.Range("G15").Value = 0.04559766666666667
$ws.Range("H15").Value = 0.136793
$ws.Range("I15").Value = 0.0005384517417354892
$ws.Range("J15").Value = 0.0006037813885615125
$ws.Range("M15").Value = 57.4814495
$ws.Range("N15").Value = 114.962899
$ws.Range("O15").Value = 0.8828264325012393
$ws.Range("P15").Value = 0.8339669888701803
$ws.Range("Q15").Value = 2.621019973817833
$ws.Range("R15").Value = 15.726119842907
$ws.Range("S15").Value = 0.0004753594302304206
$ws.Range("T15").Value = 0.0005035337465545009

$ws.Range("G16").Value = 0.04559766666666667
$ws.Range("H16").Value = 0.136793
$ws.Range("I16").Value = 0.0005384517417354892
$ws.Range("J16").Value = 0.0006037813885615125
$ws.Range("M16").Value = 1.347395
$ws.Range("N16").Value = 4.042185
$ws.Range("O16").Value = 0.02069390962418245
$ws.Range("P16").Value = 0.02932292837279799
$ws.Range("Q16").Value = 0.06143806807833333
$ws.Range("R16").Value = 0.552942612705
$ws.Range("S16").Value = 0.00001114267168045784
$ws.Range("T16").Value = 0.00001770463840961774

$ws.Range("G17").Value = 27.488287
$ws.Range("H17").Value = 54.976574
$ws.Range("I17").Value = 0.3246024872429512
$ws.Range("J17").Value = 0.2426573888143015
$ws.Range("M17").Value = 2.035934
$ws.Range("N17").Value = 6.107802
$ws.Range("O17").Value = 0.03126880699186227
$ws.Range("P17").Value = 0.04430738339814538
$ws.Range("Q17").Value = 55.96433810505799
$ws.Range("R17").Value = 335.786028630348
$ws.Range("S17").Value = 0.01014993252267827
$ws.Range("T17").Value = 0.01075151396058809

$ws.Range("G18").Value = 27.488287
$ws.Range("H18").Value = 54.976574
$ws.Range("I18").Value = 0.3246024872429512
$ws.Range("J18").Value = 0.2426573888143015
$ws.Range("O18").Value = 0.0582841555718936
$ws.Range("P18").Value = 0.08258768643246023
$ws.Range("Q18").Value = 104.3159142413763
$ws.Range("R18").Value = 625.895485448258
$ws.Range("S18").Value = 0.01891918186549177
$ws.Range("T18").Value = 0.02004051233791511

$ws.Range("G19").Value = 27.488287
$ws.Range("H19").Value = 54.976574
$ws.Range("I19").Value = 0.3246024872429512
$ws.Range("J19").Value = 0.2426573888143015
$ws.Range("M19").Value = 0.451002
$ws.Range("N19").Value = 1.353006
$ws.Range("O19").Value = 0.006926695310822388
$ws.Range("P19").Value = 0.009815012926416261
$ws.Range("Q19").Value = 12.397272413574
$ws.Range("R19").Value = 74.38363448144399
$ws.Range("S19").Value = 0.002248422526267034
$ws.Range("T19").Value = 0.002381685407902785

$ws.Range("G20").Value = 27.488287
$ws.Range("H20").Value = 54.976574
$ws.Range("I20").Value = 0.3246024872429512
$ws.Range("J20").Value = 0.2426573888143015
$ws.Range("M20").Value = 57.4814495
$ws.Range("N20").Value = 114.962899
$ws.Range("O20").Value = 0.8828264325012393
$ws.Range("P20").Value = 0.8339669888701803
$ws.Range("Q20").Value = 1580.066581032006
$ws.Range("R20").Value = 6320.266324128025
$ws.Range("S20").Value = 0.2865676557937236
$ws.Range("T20").Value = 0.2023682518765635

$ws.Range("G21").Value = 27.488287
$ws.Range("H21").Value = 54.976574
$ws.Range("I21").Value = 0.3246024872429512
$ws.Range("J21").Value = 0.2426573888143015
$ws.Range("M21").Value = 1.347395
$ws.Range("N21").Value = 4.042185
$ws.Range("O21").Value = 0.02069390962418245
$ws.Range("P21").Value = 0.02932292837279799
$ws.Range("Q21").Value = 37.037580462365
$ws.Range("R21").Value = 222.22548277419
$ws.Range("S21").Value = 0.006717294534790467
$ws.Range("T21").Value = 0.007115425231331953

$ws.Range("G22").Value = 0.027723
$ws.Range("H22").Value = 0.083169
$ws.Range("I22").Value = 0.0003273741559027063
$ws.Range("J22").Value = 0.0003670940348210247
$ws.Range("M22").Value = 2.035934
$ws.Range("N22").Value = 6.107802
$ws.Range("O22").Value = 0.03126880699186227
$ws.Range("P22").Value = 0.04430738339814538
$ws.Range("Q22").Value = 0.05644219828199999
$ws.Range("R22").Value = 0.507979784538
$ws.Range("S22").Value = 0.00001023659929504555
$ws.Range("T22").Value = 0.00001626497614398727

$ws.Range("G23").Value = 0.027723
$ws.Range("H23").Value = 0.083169
$ws.Range("I23").Value = 0.0003273741559027063
$ws.Range("J23").Value = 0.0003670940348210247
$ws.Range("O23").Value = 0.0582841555718936
$ws.Range("P23").Value = 0.08258768643246023
$ws.Range("Q23").Value = 0.105206631847
$ws.Range("R23").Value = 0.9468596866230001
$ws.Range("S23").Value = 0.00001908072623285068
$ws.Range("T23").Value = 0.00003031744703902542

$ws.Range("G24").Value = 0.027723
$ws.Range("H24").Value = 0.083169
$ws.Range("I24").Value = 0.0003273741559027063
$ws.Range("J24").Value = 0.0003670940348210247
$ws.Range("M24").Value = 0.451002
$ws.Range("N24").Value = 1.353006
$ws.Range("O24").Value = 0.006926695310822388
$ws.Range("P24").Value = 0.009815012926416261
$ws.Range("Q24").Value = 0.012503128446
$ws.Range("R24").Value = 0.112528156014
$ws.Range("S24").Value = 0.000002267621030575713
$ws.Range("T24").Value = 0.000003603032696978659

$ws.Range("G25").Value = 0.027723
$ws.Range("H25").Value = 0.083169
$ws.Range("I25").Value = 0.0003273741559027063
$ws.Range("J25").Value = 0.0003670940348210247
$ws.Range("M25").Value = 57.4814495
$ws.Range("N25").Value = 114.962899
$ws.Range("O25").Value = 0.8828264325012393
$ws.Range("P25").Value = 0.8339669888701803
$ws.Range("Q25").Value = 1.5935582244885
$ws.Range("R25").Value = 9.561349346931
$ws.Range("S25").Value = 0.0002890145581486907
$ws.Range("T25").Value = 0.0003061443068518951

$ws.Range("G26").Value = 0.027723
$ws.Range("H26").Value = 0.083169
$ws.Range("I26").Value = 0.0003273741559027063
$ws.Range("J26").Value = 0.0003670940348210247
$ws.Range("M26").Value = 1.347395
$ws.Range("N26").Value = 4.042185
$ws.Range("O26").Value = 0.02069390962418245
$ws.Range("P26").Value = 0.02932292837279799
$ws.Range("Q26").Value = 0.037353831585
$ws.Range("R26").Value = 0.336184484265
$ws.Range("S26").Value = 0.000006774651195543618
$ws.Range("T26").Value = 0.00001076427208913832

